# Generate Report for Handback
# Updates the localization-status workbook so that the handback for
# "6b3cdc19-273c-4ca3-bc19-f2813c453548.md" is reflected as completed
# ("Handed back: in sync with en-US") on all three sheets (Overview,
# zh-cn, de-de), re-sorting each table alphabetically by file name and
# refreshing the handoff/handback timestamps + clearing the stale
# "version mismatch" error message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 <- (was ffff5c51e153 row) now becomes the 6b3cdc19 row
$ws.Hyperlinks.Item(1).TextToDisplay = "e2e\6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("A2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("G2").Value = "2016-10-19 17:25:38"

# Row 3 <- (was ffffff885ab83d row) now becomes the ffff5c51e153 row
$ws.Hyperlinks.Item(2).TextToDisplay = "e2e\ffff5c51e153-17a0-45c2-a144-5ce8c92454b5.md"
$ws.Range("A3").Value = "ffff5c51e153-17a0-45c2-a144-5ce8c92454b5.md"

# Row 4 <- (was 6b3cdc19 row) now becomes the ffffff885ab83d row, handed back
$ws.Hyperlinks.Item(3).TextToDisplay = "e2e\ffffff885ab83d-8861-4568-915e-692dbfc49278.md"
$ws.Range("A4").Value = "ffffff885ab83d-8861-4568-915e-692dbfc49278.md"
$ws.Range("E4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "Handed back: in sync with en-US"
$ws.Range("G4").Value = "2016-10-19 17:21:37"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 <- becomes the 6b3cdc19 row, now handed back
$ws.Hyperlinks.Item(1).TextToDisplay = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("A2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("G2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.4357fcfaa92b9806a5ac61cede8c0228015a36ae.zh-cn.xlf"
$ws.Range("H2").Value = "2016-10-19 17:25:26"
$ws.Hyperlinks.Item(2).TextToDisplay = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("I2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("J2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.4357fcfaa92b9806a5ac61cede8c0228015a36ae.zh-cn.xlf"
$ws.Range("K2").Value = "2016-10-19 17:26:06"

# Row 3 <- becomes the ffff5c51e153 row
$ws.Hyperlinks.Item(3).TextToDisplay = "ffff5c51e153-17a0-45c2-a144-5ce8c92454b5.md"
$ws.Range("A3").Value = "ffff5c51e153-17a0-45c2-a144-5ce8c92454b5.md"
$ws.Range("F3").Value = "False"

# Row 4 <- becomes the ffffff885ab83d row, handed back
$ws.Hyperlinks.Item(5).TextToDisplay = "ffffff885ab83d-8861-4568-915e-692dbfc49278.md"
$ws.Range("A4").Value = "ffffff885ab83d-8861-4568-915e-692dbfc49278.md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "809ad35f-e201-4901-bf62-3da40928560a.df498e7df2bc578060f5d76efae911e4c169d168.zh-cn.xlf"
$ws.Range("H4").Value = "2016-10-19 17:21:26"
$ws.Hyperlinks.Item(6).TextToDisplay = "809ad35f-e201-4901-bf62-3da40928560a.md"
$ws.Range("I4").Value = "809ad35f-e201-4901-bf62-3da40928560a.md"
$ws.Range("J4").Value = "809ad35f-e201-4901-bf62-3da40928560a.df498e7df2bc578060f5d76efae911e4c169d168.zh-cn.xlf"
$ws.Range("K4").Value = "2016-10-19 17:22:09"
$ws.Range("P4").Value = ""

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 <- becomes the 6b3cdc19 row, now handed back
$ws.Hyperlinks.Item(1).TextToDisplay = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("A2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("G2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.4357fcfaa92b9806a5ac61cede8c0228015a36ae.de-de.xlf"
$ws.Range("H2").Value = "2016-10-19 17:25:38"
$ws.Hyperlinks.Item(2).TextToDisplay = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("I2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.md"
$ws.Range("J2").Value = "6b3cdc19-273c-4ca3-bc19-f2813c453548.4357fcfaa92b9806a5ac61cede8c0228015a36ae.de-de.xlf"
$ws.Range("K2").Value = "2016-10-19 17:26:24"

# Row 3 <- becomes the ffff5c51e153 row
$ws.Hyperlinks.Item(3).TextToDisplay = "ffff5c51e153-17a0-45c2-a144-5ce8c92454b5.md"
$ws.Range("A3").Value = "ffff5c51e153-17a0-45c2-a144-5ce8c92454b5.md"
$ws.Range("F3").Value = "False"

# Row 4 <- becomes the ffffff885ab83d row, handed back
$ws.Hyperlinks.Item(5).TextToDisplay = "ffffff885ab83d-8861-4568-915e-692dbfc49278.md"
$ws.Range("A4").Value = "ffffff885ab83d-8861-4568-915e-692dbfc49278.md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "809ad35f-e201-4901-bf62-3da40928560a.df498e7df2bc578060f5d76efae911e4c169d168.de-de.xlf"
$ws.Range("H4").Value = "2016-10-19 17:21:37"
$ws.Hyperlinks.Item(6).TextToDisplay = "809ad35f-e201-4901-bf62-3da40928560a.md"
$ws.Range("I4").Value = "809ad35f-e201-4901-bf62-3da40928560a.md"
$ws.Range("J4").Value = "809ad35f-e201-4901-bf62-3da40928560a.df498e7df2bc578060f5d76efae911e4c169d168.de-de.xlf"
$ws.Range("K4").Value = "2016-10-19 17:22:27"
$ws.Range("P4").Value = ""
